$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.368471145629883
$ws.Range("B1").Value = 5.728126049041748
$ws.Range("C1").Value = 6.797221660614014
$ws.Range("D1").Value = 9.444619178771973
$ws.Range("E1").Value = 5.247567653656006
